# Generate Report for Handoff
# Mark the 51463bc3-8396-4644-9193-f9b667c3f37a file as "Ready for handoff"
# across the Overview, zh-cn, and de-de report sheets, and stamp the new
# handoff date/time on each.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 51463bc3-...-f9b667c3f37a.md (row 3) ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-13-12 12:13:11"

# --- zh-cn sheet: row for 51463bc3-...-f9b667c3f37a (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-12 12:13:08"

# --- de-de sheet: row for 51463bc3-...-f9b667c3f37a (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-12 12:13:11"
